# Update "Horarios" workbook with the latest scrape for Línea 141.
# New scrape timestamp replaces the previous one across all sheets,
# two rows on the "LP1912" sheet get revised Hora_Llegada/Minutos values,
# and two brand-new rows (14_ABASTO, 81_EL PELIGRO) are appended there.

$wb = $excel.ActiveWorkbook

$oldStamp = "01:30:59"
$newStamp = "02:07:19"

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newStamp"
$ws1.Range("A3").Value = "Total filas: 4"

# Row 6: 15_ABASTO
$ws1.Cells.Item(6, 1).Value = $newStamp
$ws1.Cells.Item(6, 2).Value = "03:02"
$ws1.Cells.Item(6, 3).Value = "15_ABASTO"
$ws1.Cells.Item(6, 4).Value = 55
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Row 7: 215_ALUAR
$ws1.Cells.Item(7, 1).Value = $newStamp
$ws1.Cells.Item(7, 2).Value = "03:07"
$ws1.Cells.Item(7, 3).Value = "215_ALUAR"
$ws1.Cells.Item(7, 4).Value = 60
$ws1.Cells.Item(7, 5).Value = "LP1912"

# Row 8 (new): 14_ABASTO
$ws1.Cells.Item(8, 1).Value = $newStamp
$ws1.Cells.Item(8, 2).Value = "03:48"
$ws1.Cells.Item(8, 3).Value = "14_ABASTO"
$ws1.Cells.Item(8, 4).Value = 101
$ws1.Cells.Item(8, 5).Value = "LP1912"

# Row 9 (new): 81_EL PELIGRO
$ws1.Cells.Item(9, 1).Value = $newStamp
$ws1.Cells.Item(9, 2).Value = "04:02"
$ws1.Cells.Item(9, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(9, 4).Value = 115
$ws1.Cells.Item(9, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newStamp"

$ws2.Cells.Item(6, 1).Value = $newStamp
$ws2.Cells.Item(6, 2).Value = "03:07"
$ws2.Cells.Item(6, 3).Value = "215_ALUAR"
$ws2.Cells.Item(6, 4).Value = 60
$ws2.Cells.Item(6, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newStamp"
